# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect refreshed counts from the data source.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 40
$wsExhibit.Range("F4").Value = 16313
$wsExhibit.Range("F6").Value = 18
$wsExhibit.Range("F8").Value = 15585
$wsExhibit.Range("F10").Value = 9238
$wsExhibit.Range("F17").Value = 220
$wsExhibit.Range("F19").Value = 88
$wsExhibit.Range("F20").Value = 605
$wsExhibit.Range("F24").Value = 1150
$wsExhibit.Range("F28").Value = 522
$wsExhibit.Range("F35").Value = 265
$wsExhibit.Range("F39").Value = 5671

# Sheet "全部类型" (sheet4.xml)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 40
$wsAll.Range("F4").Value = 16313
$wsAll.Range("F6").Value = 18
$wsAll.Range("F8").Value = 15585
$wsAll.Range("F10").Value = 9238
$wsAll.Range("F17").Value = 220
$wsAll.Range("F19").Value = 88
$wsAll.Range("F20").Value = 605
$wsAll.Range("F24").Value = 1150
$wsAll.Range("F28").Value = 522
$wsAll.Range("F37").Value = 265
$wsAll.Range("F41").Value = 5671
